$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1199.4
$ws.Range("I9").Value = 483.33334
$ws.Range("J9").Value = 2273.5
$ws.Range("K9").Value = 483.33334
$ws.Range("L9").Value = 2273.5
$ws.Range("M9").Value = -314.33334
$ws.Range("N9").Value = -2611.5
$ws.Range("H32").Value = 2762
$ws.Range("I32").Value = 2762
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2762
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2436
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 583.46155
$ws.Range("I33").Value = 562.4545000000001
$ws.Range("K33").Value = 562.4545000000001
$ws.Range("M33").Value = -333.4545000000001
$ws.Range("H49").Value = 945
$ws.Range("I49").Value = 1700
$ws.Range("J49").Value = 190
$ws.Range("K49").Value = 5100
$ws.Range("L49").Value = 570
$ws.Range("M49").Value = -4964
$ws.Range("N49").Value = -842
$ws.Range("H64").Value = 5999
$ws.Range("J64").Value = 5999
$ws.Range("L64").Value = 5999
$ws.Range("N64").Value = -6495
$ws.Range("H67").Value = 5999
$ws.Range("J67").Value = 5999
$ws.Range("L67").Value = 5999
$ws.Range("N67").Value = -7715
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H100").Value = 2298.077
$ws.Range("I100").Value = 2166.111
$ws.Range("J100").Value = 2595
$ws.Range("K100").Value = 2166.111
$ws.Range("L100").Value = 2595
$ws.Range("M100").Value = -1625.111
$ws.Range("N100").Value = -3677
$ws.Range("H138").Value = 2798.75
$ws.Range("J138").Value = 3231.6667
$ws.Range("L138").Value = 9695.000100000001
$ws.Range("N138").Value = -19975.0001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 174.39131
$ws.Range("I5").Value = 179.3158
$ws.Range("K5").Value = 179.3158
$ws.Range("M5").Value = -67.3158
$ws.Range("H61").Value = 2044.5
$ws.Range("J61").Value = 2044.5
$ws.Range("L61").Value = 2044.5
$ws.Range("N61").Value = -2468.5
$ws.Range("H88").Value = 1418.0769
$ws.Range("I88").Value = 985.5714
$ws.Range("J88").Value = 1922.6666
$ws.Range("K88").Value = 985.5714
$ws.Range("L88").Value = 1922.6666
$ws.Range("M88").Value = -579.5714
$ws.Range("N88").Value = -2734.6666
$ws.Range("H91").Value = 1418.0769
$ws.Range("I91").Value = 985.5714
$ws.Range("J91").Value = 1922.6666
$ws.Range("K91").Value = 985.5714
$ws.Range("L91").Value = 1922.6666
$ws.Range("M91").Value = 418.4286
$ws.Range("N91").Value = -4730.6666
$ws.Range("H102").Value = 2056.5833
$ws.Range("I102").Value = 1653.2222
$ws.Range("K102").Value = 1653.2222
$ws.Range("M102").Value = -31.22219999999993
$ws.Range("H132").Value = 2901.6
$ws.Range("I132").Value = 2630.75
$ws.Range("K132").Value = 7892.25
$ws.Range("M132").Value = -5362.25
$ws.Range("H136").Value = 2044.5
$ws.Range("J136").Value = 2044.5
$ws.Range("L136").Value = 6133.5
$ws.Range("N136").Value = -11233.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 174.39131
$ws.Range("I4").Value = 179.3158
$ws.Range("K4").Value = 179.3158
$ws.Range("M4").Value = -64.3158
$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -21058
$ws.Range("H33").Value = 80000
$ws.Range("J33").Value = 80000
$ws.Range("L33").Value = 80000
$ws.Range("N33").Value = -80672
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H105").Value = 810.46155
$ws.Range("I105").Value = 710.4167
$ws.Range("K105").Value = 710.4167
$ws.Range("M105").Value = 1036.5833
$ws.Range("H107").Value = 2500
$ws.Range("I107").Value = 2166.6667
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 2166.6667
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = -246.6667000000002
$ws.Range("N107").Value = -7340
$ws.Range("H134").Value = 2096.8
$ws.Range("I134").Value = 1999.6666
$ws.Range("J134").Value = 2242.5
$ws.Range("K134").Value = 5998.9998
$ws.Range("L134").Value = 6727.5
$ws.Range("M134").Value = -3463.9998
$ws.Range("N134").Value = -11797.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1659.907
$ws.Range("I31").Value = 907.12823
$ws.Range("K31").Value = 907.12823
$ws.Range("M31").Value = -612.12823
$ws.Range("H34").Value = 1659.907
$ws.Range("I34").Value = 907.12823
$ws.Range("K34").Value = 907.12823
$ws.Range("M34").Value = -705.12823
$ws.Range("H44").Value = 21687.334
$ws.Range("I44").Value = 3064
$ws.Range("J44").Value = 30999
$ws.Range("K44").Value = 3064
$ws.Range("L44").Value = 30999
$ws.Range("M44").Value = -2622
$ws.Range("N44").Value = -31883
$ws.Range("H59").Value = 63737
$ws.Range("J59").Value = 63737
$ws.Range("L59").Value = 63737
$ws.Range("N59").Value = -66027
$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 229
$ws.Range("I25").Value = 211
$ws.Range("J25").Value = 289
$ws.Range("K25").Value = 633
$ws.Range("L25").Value = 867
$ws.Range("M25").Value = -464
$ws.Range("N25").Value = -1205
$ws.Range("H29").Value = 126.666664
$ws.Range("I29").Value = 31
$ws.Range("J29").Value = 174.5
$ws.Range("K29").Value = 93
$ws.Range("L29").Value = 523.5
$ws.Range("M29").Value = 184
$ws.Range("N29").Value = -1077.5
$ws.Range("H30").Value = 229
$ws.Range("I30").Value = 211
$ws.Range("J30").Value = 289
$ws.Range("K30").Value = 633
$ws.Range("L30").Value = 867
$ws.Range("M30").Value = -531
$ws.Range("N30").Value = -1071
$ws.Range("H44").Value = 906.2
$ws.Range("I44").Value = 270.5
$ws.Range("J44").Value = 1004
$ws.Range("K44").Value = 811.5
$ws.Range("L44").Value = 3012
$ws.Range("M44").Value = -413.5
$ws.Range("N44").Value = -3808
$ws.Range("H45").Value = 2749.5
$ws.Range("J45").Value = 2749.5
$ws.Range("L45").Value = 8248.5
$ws.Range("N45").Value = -9312.5
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H123").Value = 1990
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 2500
$ws.Range("I124").Value = 2500
$ws.Range("K124").Value = 7500
$ws.Range("M124").Value = -2590
$ws.Range("H139").Value = 47304
$ws.Range("I139").Value = 3390.6667
$ws.Range("K139").Value = 10172.0001
$ws.Range("M139").Value = -5032.000100000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 9031.75
$ws.Range("I22").Value = 1063.5
$ws.Range("J22").Value = 17000
$ws.Range("K22").Value = 1063.5
$ws.Range("L22").Value = 17000
$ws.Range("M22").Value = -534.5
$ws.Range("N22").Value = -18058
$ws.Range("H31").Value = 235.16667
$ws.Range("I31").Value = 235.16667
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 235.16667
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 56.83332999999999
$ws.Range("N31").ClearContents()
$ws.Range("H37").Value = 235.16667
$ws.Range("I37").Value = 235.16667
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 235.16667
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 41.83332999999999
$ws.Range("N37").ClearContents()
$ws.Range("H126").Value = 3540.1667
$ws.Range("I126").Value = 3688.2
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 11064.6
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -8594.599999999999
$ws.Range("N126").Value = -13340

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 846.25
$ws.Range("I55").Value = 1205
$ws.Range("J55").Value = 487.5
$ws.Range("K55").Value = 1205
$ws.Range("L55").Value = 487.5
$ws.Range("M55").Value = -1032
$ws.Range("N55").Value = -833.5
$ws.Range("H104").Value = 27346.223
$ws.Range("J104").Value = 27346.223
$ws.Range("L104").Value = 27346.223
$ws.Range("N104").Value = -34334.223

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5999.6665
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 5999.5
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 5999.5
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -7247.5
$ws.Range("H65").Value = 5999.6665
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 5999.5
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 29997.5
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -36237.5
$ws.Range("H68").Value = 34408.4
$ws.Range("J68").Value = 34408.4
$ws.Range("L68").Value = 34408.4
$ws.Range("N68").Value = -36030.4
$ws.Range("H71").Value = 34408.4
$ws.Range("J71").Value = 34408.4
$ws.Range("L71").Value = 103225.2
$ws.Range("N71").Value = -111337.2
